# Update countries & provincias Spain
# Refreshes the COVID-19 country table ("Pais" sheet) with newer totals and
# updates the "last updated" timestamp. Because the sheet is kept sorted
# descending by "Casos totales" (column B), a handful of countries swapped
# rank with their neighbour, so those row pairs are re-written together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    if ($country -ne $null) {
        $ws.Cells.Item($row, 1).Value = $country
    }
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Timestamp footer (row 1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 17:57"

# --- Straight numeric refreshes (country stays in the same row) -----------
Set-Row 4   $null 5427542 11876 2844525 2412285 0 317 170732
Set-Row 6   $null 2506247 46634 1770682  686677 0 744  48888
Set-Row 15  $null  313798     0       0       0 0  11  41358
Set-Row 20  $null  252809   574  203326   14249 0   3  35234
Set-Row 47  $null   55580    83   51049    4504 0   0     27
Set-Row 49  $null   53783   235   39374   12637 0   2   1772
Set-Row 58  $null   37671   268   32900    2780 0   0   1991
Set-Row 102 $null    6632   251    3804    2605 0   2    223

# --- Libia / Albania swap places (rows 99-100) -----------------------------
Set-Row 99  "Albania" 7117 146 3695 3203 0 6 219
Set-Row 100 "Libia"   7050   0  816 6099 0 0 135

# --- Republica de Chipre jumps ahead of Letonia and Georgia (rows 145-147) -
Set-Row 145 "Republica de Chipre" 1318 13  870  428 0 0 20
Set-Row 146 "Letonia"             1308  1 1078  198 0 0 32
Set-Row 147 "Georgia"             1306 23 1085  204 0 0 17

# --- Trinidad yTobago / Burundi swap places (rows 168-169) -----------------
Set-Row 168 "Trinidad yTobago" 412 8 139 265 0 0 8
Set-Row 169 "Burundi"          410 0 315  94 0 0 1

# --- Montserrat / Islas Malvinas swap places (rows 213-214) ----------------
Set-Row 213 "Montserrat"      13 0 12 0 0 0 1
Set-Row 214 "Islas Malvinas"  13 0 13 0 0 0 0
